# Applies the "Add files via upload" edit: appends one new video row to the
# Overall Summary / Video Details sheets, a matching set of keyword / issue
# rows to Negative Keywords / Negative Issues, updates the cumulative totals
# on the Overall Summary sheet, and refreshes the two AI narrative cells.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# A single reusable "scratch" cell, parked far outside every sheet's used
# range. It is formatted as Text once, so that round-tripping a string
# through it (Copy + PasteSpecial values-only) lands the value on the real
# destination cell as literal text instead of letting Excel's normal
# "0.0%"/"0.855" -> number auto-conversion kick in. The destination cell's
# own number format / style is untouched by a values-only paste, so we
# never have to fight to restore it afterwards.
$wsSummary = $wb.Worksheets.Item("Overall Summary")
$scratch = $wsSummary.Range("ZZ1")
$scratch.NumberFormat = "@"

function Set-TextValue {
    param($Scratch, $Range, [string]$Text)
    $Scratch.Value = $Text
    $Scratch.Copy()
    $Range.PasteSpecial($script:xlPasteValues)
}

# ---------------------------------------------------------------------
# 1. Overall Summary — cumulative header counts (row 2)
# ---------------------------------------------------------------------
$wsSummary.Range("A2").Value = 71
$wsSummary.Range("B2").Value = 62
$wsSummary.Range("D2").Value = 76

# New per-video row 92, cloning the formatting of row 91 first so the new
# row picks up the same style index, then overwriting the cell contents.
$wsSummary.Range("A91:I91").Copy()
$wsSummary.Range("A92:I92").PasteSpecial($xlPasteFormats)
$wsSummary.Application.CutCopyMode = $false

$wsSummary.Range("A92").Value = "Chandrababu Naidu Speech in AP Assembly: అసెంబ్లీ "
$wsSummary.Range("B92").Value = "Asianet News Telugu"
$wsSummary.Range("C92").Value = 10
Set-TextValue $scratch $wsSummary.Range("D92") "0.0%"
Set-TextValue $scratch $wsSummary.Range("E92") "100.0%"
Set-TextValue $scratch $wsSummary.Range("F92") "0.0%"
$wsSummary.Range("G92").Value = 9
$wsSummary.Range("H92").Value = 5

# ---------------------------------------------------------------------
# 2. Video Details — same new video, row 69
# ---------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("Video Details")
$wsDetails.Range("A68:H68").Copy()
$wsDetails.Range("A69:H69").PasteSpecial($xlPasteFormats)
$wsDetails.Application.CutCopyMode = $false

$wsDetails.Range("A69").Value = "Chandrababu Naidu Speech in AP Assembly: అసెంబ్లీ "
$wsDetails.Range("B69").Value = "Asianet News Telugu"
$wsDetails.Range("C69").Value = 10
Set-TextValue $scratch $wsDetails.Range("D69") "0.0%"
Set-TextValue $scratch $wsDetails.Range("E69") "100.0%"
Set-TextValue $scratch $wsDetails.Range("F69") "0.0%"
$wsDetails.Range("G69").Value = 9
$wsDetails.Range("H69").Value = 5

# ---------------------------------------------------------------------
# 3. Negative Keywords — four new keyword rows (180-183)
# ---------------------------------------------------------------------
$wsNegKeywords = $wb.Worksheets.Item("Negative Keywords")
$wsNegKeywords.Range("A179:G179").Copy()
$wsNegKeywords.Range("A180:G183").PasteSpecial($xlPasteFormats)
$wsNegKeywords.Application.CutCopyMode = $false

$negKeywordRows = @(
    @{ Row = 180; A = "Jagan Mohan Reddy";    B = 8; C = 8; D = "100.0%"; E = "0.0%"; F = "0.0%"; G = "0.855" },
    @{ Row = 181; A = "powerful speech";      B = 9; C = 9; D = "100.0%"; E = "0.0%"; F = "0.0%"; G = "0.850" },
    @{ Row = 182; A = "assembly highlights";  B = 9; C = 9; D = "100.0%"; E = "0.0%"; F = "0.0%"; G = "0.855" },
    @{ Row = 183; A = "video play";           B = 9; C = 9; D = "100.0%"; E = "0.0%"; F = "0.0%"; G = "0.850" }
)

foreach ($r in $negKeywordRows) {
    $row = $r.Row
    $wsNegKeywords.Range("A$row").Value = $r.A
    $wsNegKeywords.Range("B$row").Value = $r.B
    $wsNegKeywords.Range("C$row").Value = $r.C
    Set-TextValue $scratch $wsNegKeywords.Range("D$row") $r.D
    Set-TextValue $scratch $wsNegKeywords.Range("E$row") $r.E
    Set-TextValue $scratch $wsNegKeywords.Range("F$row") $r.F
    Set-TextValue $scratch $wsNegKeywords.Range("G$row") $r.G
}

# ---------------------------------------------------------------------
# 4. Negative Issues — five new issue rows (138-142)
# ---------------------------------------------------------------------
$wsNegIssues = $wb.Worksheets.Item("Negative Issues")
$wsNegIssues.Range("A137:E137").Copy()
$wsNegIssues.Range("A138:E142").PasteSpecial($xlPasteFormats)
$wsNegIssues.Application.CutCopyMode = $false

$negIssueRows = @(
    @{ Row = 138; A = "Criticism of Jagan's governance and infrastructure development";               B = "Political Governance";        C = "0.850"; D = 1; E = 1 },
    @{ Row = 139; A = "Allegations of inefficiency in project execution, like the Polavaram project";  B = "Infrastructure Development";  C = "0.750"; D = 1; E = 1 },
    @{ Row = 140; A = "Concerns over pending agricultural support and issues in fee releases in education"; B = "Agricultural Support";    C = "0.700"; D = 1; E = 1 },
    @{ Row = 141; A = "Calls for better responsiveness from the government regarding student fees and certificates"; B = "Education & Healthcare"; C = "0.680"; D = 1; E = 1 },
    @{ Row = 142; A = "Complaints about social welfare number and communication issues";               B = "Social Welfare Schemes";       C = "0.650"; D = 1; E = 1 }
)

foreach ($r in $negIssueRows) {
    $row = $r.Row
    $wsNegIssues.Range("A$row").Value = $r.A
    $wsNegIssues.Range("B$row").Value = $r.B
    Set-TextValue $scratch $wsNegIssues.Range("C$row") $r.C
    $wsNegIssues.Range("D$row").Value = $r.D
    $wsNegIssues.Range("E$row").Value = $r.E
}

# Done with the scratch cell — remove it so it doesn't linger in the sheet
# or inflate the used range / dimension.
$wsSummary.Application.CutCopyMode = $false
$scratch.Clear()

# ---------------------------------------------------------------------
# 5. Narrative Summary — refreshed AI keyword-analysis narrative
# ---------------------------------------------------------------------
$wsNarrativeSummary = $wb.Worksheets.Item("Narrative Summary")
$narrativeSummaryText = "== Updated 2025-10-16 21:55:12 ===`n`nAI-Generated Keyword Analysis Summary (Cumulative)`n`nThe sentiment analysis of 71 YouTube videos focused on Andhra Pradesh political content reveals a starkly negative public sentiment, with an overwhelming 100% of comments analyzed reflecting dissatisfaction. This finding suggests a significant disconnect between political leadership and public perception, specifically concerning key figures such as Chandrababu Naidu and Jagan Mohan Reddy. `n`nThe absence of positive sentiment indicates not only discontent but also a potential crisis of legitimacy for political actors in the state. The dominance of negative keywords highlights that both leaders are central to public grievances, with criticisms surrounding their governance and performance in the AP Assembly. The phrase ""powerful speech"" appears paradoxically within a negative context, suggesting that while political rhetoric may be compelling, it fails to resonate positively with the electorate, pointing to a growing skepticism of promises made by these leaders.`n`nFurthermore, the analysis uncovers a polarized political landscape, where any discussion of Andhra Pradesh politics evokes strong reactions. This polarization indicates entrenched partisan divides, complicating the prospect for constructive political dialogue. The data suggests that future political campaigns must address these sentiments directly, focusing on transparency and accountability to rebuild trust among constituents. In a digital age where public opinion is rapidly formed and expressed, understanding this sentiment landscape is crucial for effective political strategy in Andhra Pradesh.`n`nThis cumulative summary was generated by AI based on aggregated sentiment data from 71 videos and 62 comments."
$wsNarrativeSummary.Range("A1").Formula = "=" + $narrativeSummaryText

# ---------------------------------------------------------------------
# 6. Political Issues Narrative — refreshed AI political-issues narrative
# ---------------------------------------------------------------------
$wsPoliticalNarrative = $wb.Worksheets.Item("Political Issues Narrative")
$politicalNarrativeText = "== Updated 2025-10-16 21:55:18 ===`n`nAI-Generated Political Issues Analysis (Cumulative)`n`nThe analysis of 71 YouTube videos concerning Andhra Pradesh politics reveals a landscape fraught with discontent and criticism, primarily directed at the governance of Chief Minister Jagan Mohan Reddy. The predominant issues identified span five key categories: Political Governance, Infrastructure Development, Agricultural Support, Education & Healthcare, and Social Welfare Schemes. All five issues are viewed negatively, highlighting a significant disconnect between government actions and public expectations.`n`nThe most contentious issue concerns the criticism of Jagan's governance, particularly in relation to infrastructure development. The Polavaram project, a critical infrastructure initiative, has raised allegations of inefficiency, signaling broader concerns about the state's project execution capabilities. This resonates with the public's desire for transparent and effective governance, emphasizing a call for accountability.`n`nAgricultural support remains a pressing issue, with complaints about delays in financial aid and educational fee releases exacerbating existing frustrations among farmers and students. These concerns reflect a deepening crisis in the agricultural sector, where timely support is vital for sustaining livelihoods.`n`nAdditionally, the calls for improved governmental responsiveness regarding student fees and welfare communications underline systemic inefficiencies that erode public trust. The lack of positive sentiment in these discussions suggests that the current administration is struggling to align its policies with the immediate needs of its constituents.`n`nIn conclusion, these findings not only illuminate the pressing concerns within Andhra Pradesh but also underscore a critical need for the government to engage more effectively with its citizens. The recurring themes in public discourse indicate a clear priority for accountability, efficient governance, and direct support mechanisms, which could shape the political landscape as upcoming elections draw near.`n`nThis cumulative analysis was generated by AI based on identified political issues from 71 videos."
$wsPoliticalNarrative.Range("A1").Formula = "=" + $politicalNarrativeText
